$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 547.3913
$ws.Range("I28").Value = 374.7143
$ws.Range("J28").Value = 2360.5
$ws.Range("K28").Value = 374.7143
$ws.Range("L28").Value = 2360.5
$ws.Range("M28").Value = 110.2857
$ws.Range("N28").Value = -3330.5

$ws.Range("H33").Value = 184.5
$ws.Range("I33").Value = 151.21053
$ws.Range("K33").Value = 151.21053
$ws.Range("M33").Value = 77.78946999999999

$ws.Range("H38").Value = 72.625
$ws.Range("I38").Value = 72.625
$ws.Range("K38").Value = 217.875
$ws.Range("M38").Value = 154.125

$ws.Range("H113").Value = 2633.4666
$ws.Range("I113").Value = 2333.5
$ws.Range("K113").Value = 2333.5
$ws.Range("M113").Value = 920.5

$ws.Range("H132").Value = 9149.25
$ws.Range("I132").Value = 9414.037
$ws.Range("K132").Value = 28242.111
$ws.Range("M132").Value = -25712.111

$ws.Range("H135").Value = 2935.2727
$ws.Range("I135").Value = 1143.2222
$ws.Range("J135").Value = 10999.5
$ws.Range("K135").Value = 10288.9998
$ws.Range("L135").Value = 98995.5
$ws.Range("M135").Value = -7753.9998
$ws.Range("N135").Value = -104065.5

$ws.Range("H137").Value = 10646133
$ws.Range("I137").Value = 20002570
$ws.Range("J137").Value = 13819.5
$ws.Range("K137").Value = 60007710
$ws.Range("L137").Value = 41458.5
$ws.Range("M137").Value = -60005160
$ws.Range("N137").Value = -46558.5

$ws.Range("H138").Value = 8376.166999999999
$ws.Range("I138").Value = 11619.8
$ws.Range("J138").Value = 7937.838
$ws.Range("K138").Value = 34859.39999999999
$ws.Range("L138").Value = 23813.514
$ws.Range("M138").Value = -29719.39999999999
$ws.Range("N138").Value = -34093.514

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 794716
$ws.Range("I32").Value = 895441.0600000001
$ws.Range("J32").Value = 14096.75
$ws.Range("K32").Value = 895441.0600000001
$ws.Range("L32").Value = 14096.75
$ws.Range("M32").Value = -895154.0600000001
$ws.Range("N32").Value = -14670.75

$ws.Range("H45").Value = 7974.25
$ws.Range("I45").Value = 7974.25
$ws.Range("K45").Value = 7974.25
$ws.Range("M45").Value = -7597.25

$ws.Range("H61").Value = 2328800.8
$ws.Range("I61").Value = 2959.3928
$ws.Range("K61").Value = 2959.3928
$ws.Range("M61").Value = -2747.3928

$ws.Range("H74").Value = 486714.62
$ws.Range("I74").Value = 545725.4399999999
$ws.Range("J74").Value = 14628.2
$ws.Range("K74").Value = 545725.4399999999
$ws.Range("L74").Value = 14628.2
$ws.Range("M74").Value = -544851.4399999999
$ws.Range("N74").Value = -16376.2

$ws.Range("H77").Value = 486714.62
$ws.Range("I77").Value = 545725.4399999999
$ws.Range("J77").Value = 14628.2
$ws.Range("K77").Value = 2728627.2
$ws.Range("L77").Value = 73141
$ws.Range("M77").Value = -2724259.2
$ws.Range("N77").Value = -81877

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 1434.7391
$ws.Range("I122").Value = 1257.2222
$ws.Range("J122").Value = 2073.8
$ws.Range("K122").Value = 3771.6666
$ws.Range("L122").Value = 6221.400000000001
$ws.Range("M122").Value = -1321.6666
$ws.Range("N122").Value = -11121.4

$ws.Range("H132").Value = 4610.4365
$ws.Range("I132").Value = 3279.9285
$ws.Range("J132").Value = 5990.222
$ws.Range("K132").Value = 9839.7855
$ws.Range("L132").Value = 17970.666
$ws.Range("M132").Value = -7309.7855
$ws.Range("N132").Value = -23030.666

$ws.Range("H136").Value = 2328800.8
$ws.Range("I136").Value = 2959.3928
$ws.Range("K136").Value = 8878.178400000001
$ws.Range("M136").Value = -6328.178400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3012.0833
$ws.Range("I86").Value = 3005.2222
$ws.Range("K86").Value = 3005.2222
$ws.Range("M86").Value = -1882.2222

$ws.Range("H89").Value = 3012.0833
$ws.Range("I89").Value = 3005.2222
$ws.Range("K89").Value = 15026.111
$ws.Range("M89").Value = -9410.111000000001

$ws.Range("H99").Value = 6263.8
$ws.Range("I99").Value = 6388.375
$ws.Range("J99").Value = 5765.5
$ws.Range("K99").Value = 6388.375
$ws.Range("L99").Value = 5765.5
$ws.Range("M99").Value = -4890.375
$ws.Range("N99").Value = -8761.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 12424.556
$ws.Range("I94").Value = 25552
$ws.Range("J94").Value = 1922.6
$ws.Range("K94").Value = 25552
$ws.Range("L94").Value = 1922.6
$ws.Range("M94").Value = -25101
$ws.Range("N94").Value = -2824.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 888.7
$ws.Range("I18").Value = 462
$ws.Range("J18").Value = 1173.1666
$ws.Range("K18").Value = 1386
$ws.Range("L18").Value = 3519.4998
$ws.Range("M18").Value = -1217
$ws.Range("N18").Value = -3857.4998

$ws.Range("H33").Value = 4959874.5
$ws.Range("I33").Value = 8678015
$ws.Range("J33").Value = 2354.1667
$ws.Range("K33").Value = 52068090
$ws.Range("L33").Value = 14125.0002
$ws.Range("M33").Value = -52067807
$ws.Range("N33").Value = -14691.0002

$ws.Range("H68").Value = 4923.636
$ws.Range("I68").Value = 3083.3333
$ws.Range("J68").Value = 5011.27
$ws.Range("K68").Value = 9249.999899999999
$ws.Range("L68").Value = 15033.81
$ws.Range("M68").Value = -8438.999899999999
$ws.Range("N68").Value = -16655.81

$ws.Range("H71").Value = 4923.636
$ws.Range("I71").Value = 3083.3333
$ws.Range("J71").Value = 5011.27
$ws.Range("K71").Value = 27749.9997
$ws.Range("L71").Value = 45101.43000000001
$ws.Range("M71").Value = -23693.9997
$ws.Range("N71").Value = -53213.43000000001

$ws.Range("H86").Value = 716.6667
$ws.Range("I86").Value = 700
$ws.Range("K86").Value = 2100
$ws.Range("M86").Value = -914

$ws.Range("H89").Value = 716.6667
$ws.Range("I89").Value = 700
$ws.Range("K89").Value = 6300
$ws.Range("M89").Value = -372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1603
$ws.Range("I113").Value = 1603
$ws.Range("K113").Value = 1603
$ws.Range("M113").Value = 567

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5444.909
$ws.Range("I40").Value = 4761.5
$ws.Range("J40").Value = 6265
$ws.Range("K40").Value = 4761.5
$ws.Range("L40").Value = 6265
$ws.Range("M40").Value = -4625.5
$ws.Range("N40").Value = -6537

$ws.Range("H132").Value = 2528505.2
$ws.Range("I132").Value = 2781092.5
$ws.Range("J132").Value = 2633.3333
$ws.Range("K132").Value = 8343277.5
$ws.Range("L132").Value = 7899.999899999999
$ws.Range("M132").Value = -8340747.5
$ws.Range("N132").Value = -12959.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1159.0952
$ws.Range("I113").Value = 711.9167
$ws.Range("K113").Value = 2135.7501
$ws.Range("M113").Value = 34.2498999999998

$ws.Range("H127").Value = 69980
$ws.Range("J127").Value = 69980
$ws.Range("L127").Value = 69980
$ws.Range("N127").Value = -79900
